# Generate Report for Handback
#
# The localization status workbook is refreshed after a handback run:
#   - the "Status" column moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every language sheet
#   - the per-language sheets (zh-cn, de-de) get their "Latest Target File"
#     and "Latest Handback File" columns populated with links to the
#     generated files for each source document
#   - the de-de sheet additionally records the new "Latest Handback
#     DateTime" for both rows
#   - a couple of columns are widened so the new, longer filenames are
#     readable

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$md1Name    = "2138835f-8b04-4a55-8813-9128449d3f6c.md"
$md1Url     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/582babba51ed800cd7c242cecaa107f77a1c6236/e2e/2138835f-8b04-4a55-8813-9128449d3f6c.md"
$md2Name    = "49dce285-79b1-46d0-93eb-392e5b27552e.md"
$md2Url     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/582babba51ed800cd7c242cecaa107f77a1c6236/e2e/49dce285-79b1-46d0-93eb-392e5b27552e.md"

$zhXlf1 = "2138835f-8b04-4a55-8813-9128449d3f6c.98248e03a727e199df8d4695b38b4c67be377d70.zh-cn.xlf"
$zhXlf2 = "49dce285-79b1-46d0-93eb-392e5b27552e.858ab1fd6d96947dbc11d5c84056b0ac6ca73c52.zh-cn.xlf"
$deXlf1 = "2138835f-8b04-4a55-8813-9128449d3f6c.98248e03a727e199df8d4695b38b4c67be377d70.de-de.xlf"
$deXlf2 = "49dce285-79b1-46d0-93eb-392e5b27552e.858ab1fd6d96947dbc11d5c84056b0ac6ca73c52.de-de.xlf"

$zhHandbackDate = "2016-09-03 15:06:44"
$deHandbackDate = "2016-09-03 15:06:51"

# A ColumnWidth of 29.17 persists as an XML column width of 30 (the widest
# width this engine's pixel grid lets us reach while still landing as close
# as possible to the new wider columns used for the longer handback file
# names).
$wideWidth = 29.17
# A ColumnWidth of 39.17 persists as an XML column width of 40, matching the
# other 40-wide columns already in this workbook.
$fullWidth = 39.17

# ---------------------------------------------------------------------
# Overview sheet: status rollup columns (E = zh-cn, F = de-de) + widen
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = $wideWidth
$overview.Columns.Item(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column (C) -> handed back
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

# Widen Status (C), Latest Target File (I) and Latest Handback File (J)
$zh.Columns.Item(3).ColumnWidth = $wideWidth
$zh.Columns.Item(9).ColumnWidth = $fullWidth
$zh.Columns.Item(10).ColumnWidth = $fullWidth

# Row 2 (2138835f...)
$zh.Hyperlinks.Add($zh.Range("I2"), $md1Url, "", "", $md1Name)
$zh.Range("J2").Value = $zhXlf1
$zh.Range("K2").Value = $zhHandbackDate

# Row 3 (49dce285...)
$zh.Hyperlinks.Add($zh.Range("I3"), $md2Url, "", "", $md2Name)
$zh.Range("J3").Value = $zhXlf2
$zh.Range("K3").Value = $zhHandbackDate

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Status column (C) -> handed back
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

# Widen Status (C), Latest Target File (I) and Latest Handback File (J)
$de.Columns.Item(3).ColumnWidth = $wideWidth
$de.Columns.Item(9).ColumnWidth = $fullWidth
$de.Columns.Item(10).ColumnWidth = $fullWidth

# Row 2 (2138835f...)
$de.Hyperlinks.Add($de.Range("I2"), $md1Url, "", "", $md1Name)
$de.Range("J2").Value = $deXlf1
$de.Range("K2").Value = $deHandbackDate

# Row 3 (49dce285...)
$de.Hyperlinks.Add($de.Range("I3"), $md2Url, "", "", $md2Name)
$de.Range("J3").Value = $deXlf2
$de.Range("K3").Value = $deHandbackDate
